# Updated symbol list on Wed Feb  8 19:46:55 UTC 2023 with GitHub Actions
# Applies refreshed Price (column D) / Volume(1h) (column E) values for the
# cryptos sheet. All cells on this sheet are stored as text, so we force the
# Text number format before writing to avoid Excel re-interpreting values
# like "0.02740" or "0.003490" as numbers and dropping trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "326.92";       E = "-0.73%" },
    @{ Row = 3;  D = "44.45";        E = "0.23%" },
    @{ Row = 4;  D = "5.245";        E = "-4.95%" },
    @{ Row = 5;  D = "0.08373";      E = "3.39%" },
    @{ Row = 6;  D = "1.939";        E = "-5.87%" },
    @{ Row = 7;  D = "0.9719";       E = "-0.17%" },
    @{ Row = 8;  D = "2.536";        E = "-4.72%" },
    @{ Row = 9;  D = "0.1156";       E = "2.90%" },
    @{ Row = 10; D = "0.1896";       E = "0.29%" },
    @{ Row = 11; D = "0.09643";      E = "-3.39%" },
    @{ Row = 12; D = "0.04616";      E = "-3.32%" },
    @{ Row = 13;               E = "0.20%" },
    @{ Row = 14; D = "0.001293";     E = "1.70%" },
    @{ Row = 15; D = "0.005811";     E = "-3.12%" },
    @{ Row = 16; D = "3.401";        E = "1.82%" },
    @{ Row = 17; D = "4.449";        E = "0.50%" },
    @{ Row = 18;               E = "1.70%" },
    @{ Row = 19; D = "8.654";        E = "-14.99%" },
    @{ Row = 20; D = "0.1361";       E = "-2.07%" },
    @{ Row = 22; D = "0.04165";      E = "1.64%" },
    @{ Row = 23; D = "0.001235";     E = "-5.29%" },
    @{ Row = 24; D = "0.004425";     E = "0.72%" },
    @{ Row = 25; D = "0.0001304";    E = "1.95%" },
    @{ Row = 26; D = "0.0002983";    E = "-20.15%" },
    @{ Row = 38; D = "0.02740";      E = "2.34%" },
    @{ Row = 39; D = "0.05638";      E = "-0.19%" },
    @{ Row = 40; D = "0.007858";     E = "3.01%" },
    @{ Row = 41; D = "0.1409";       E = "-0.42%" },
    @{ Row = 42; D = "0.007354";     E = "-1.08%" },
    @{ Row = 43; D = "0.002122";     E = "8.39%" },
    @{ Row = 44; D = "0.007873";     E = "-5.03%" },
    @{ Row = 45; D = "0.3501" },
    @{ Row = 46; D = "0.00006862";   E = "-3.06%" },
    @{ Row = 47; D = "0.00000000751"; E = "0.24%" },
    @{ Row = 48; D = "0.003490";     E = "-1.59%" },
    @{ Row = 49; D = "0.003535";     E = "40.42%" },
    @{ Row = 50; D = "0.00002103";   E = "0.24%" },
    @{ Row = 51; D = "0.0002003";    E = "0.24%" }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.ContainsKey("D")) {
        $cellD = $ws.Cells.Item($row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
    }

    if ($u.ContainsKey("E")) {
        $cellE = $ws.Cells.Item($row, 5)
        $cellE.NumberFormat = "@"
        $cellE.Value = $u.E
    }
}
